# Update the "repaymentstrategy" value on the ProductLoanInput sheet (row 17,
# column B) from "Mifos style" to "Penalties, Fees, Interest, Principal order",
# which also introduces a new cell style (left/top aligned, no wrap) and moves
# the sheet's active selection to the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$cell = $ws.Range("B17")
$cell.Value2 = "Penalties, Fees, Interest, Principal order"

# New style: fillId=2 (green), no font override, horizontal=left, vertical=top,
# no wrap text (matches cellXfs index 14 introduced by the edit).
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160
$cell.WrapText = $false

# Move the active selection to the edited cell.
$ws.Activate() | Out-Null
$cell.Select() | Out-Null
